# Edit: fill in the title-slide placeholders with text and remove the
# second (blank "Hello" / trial) slide from the deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title placeholder ("Title 1") -> "This is a slide" (two runs)
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "This is "
$titleRange.LanguageID = "en-IN"
[void]$titleRange.InsertAfter("a slide")

# Subtitle placeholder ("Subtitle 2") -> "Hi this is trial" (two runs)
$subRange = $s.Shapes.Item(2).TextFrame.TextRange
$subRange.Text = "Hi this "
$subRange.LanguageID = "en-IN"
[void]$subRange.InsertAfter("is trial")

# Remove the second slide from the presentation entirely.
$p.Slides.Item(2).Delete()
